$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Update existing data rows (2-21): Fecha, Volumen, Precio min/max,
#    Precio promedio ponderado and Precio $/Kg were revised for the
#    weekly refresh of the data set.
# -----------------------------------------------------------------
$ws.Range("D2").Value = 44350
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19500
$ws.Range("S2").Value = 975

$ws.Range("D3").Value = 44326
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 19500
$ws.Range("P3").Value = 19750
$ws.Range("S3").Value = 988

$ws.Range("D4").Value = 44364
$ws.Range("M4").Value = 140

$ws.Range("D5").Value = 44336
$ws.Range("N5").Value = 19500
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 19750
$ws.Range("S5").Value = 988

$ws.Range("D6").Value = 44445
$ws.Range("M6").Value = 160

$ws.Range("D7").Value = 44333
$ws.Range("M7").Value = 100
$ws.Range("N7").Value = 19500
$ws.Range("O7").Value = 20000
$ws.Range("P7").Value = 19750
$ws.Range("S7").Value = 988

$ws.Range("D8").Value = 44431
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 21000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21500
$ws.Range("S8").Value = 1075

$ws.Range("D9").Value = 44418
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 21000
$ws.Range("P9").Value = 20500
$ws.Range("S9").Value = 1025

$ws.Range("D10").Value = 44407
$ws.Range("M10").Value = 160

$ws.Range("D11").Value = 44410
$ws.Range("M11").Value = 200
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20500
$ws.Range("S11").Value = 1025

$ws.Range("D12").Value = 44343
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 19500
$ws.Range("O12").Value = 20000
$ws.Range("P12").Value = 19750
$ws.Range("S12").Value = 988

$ws.Range("D13").Value = 44365
$ws.Range("M13").Value = 100
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 21000
$ws.Range("P13").Value = 20500
$ws.Range("S13").Value = 1025

$ws.Range("D14").Value = 44335
$ws.Range("M14").Value = 200
$ws.Range("N14").Value = 19000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 19500
$ws.Range("S14").Value = 975

$ws.Range("D15").Value = 44434
$ws.Range("M15").Value = 100

$ws.Range("D17").Value = 44427
$ws.Range("M17").Value = 200
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 21000
$ws.Range("P17").Value = 20500
$ws.Range("S17").Value = 1025

$ws.Range("D19").Value = 44417
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 21000
$ws.Range("P19").Value = 20500
$ws.Range("S19").Value = 1025

$ws.Range("D20").Value = 44441
$ws.Range("M20").Value = 160
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 21000
$ws.Range("P20").Value = 20500
$ws.Range("S20").Value = 1025

$ws.Range("D21").Value = 44315
$ws.Range("M21").Value = 100

# -----------------------------------------------------------------
# 2) Append three new weekly records (rows 22-24) with the same
#    shape as the existing data rows.
# -----------------------------------------------------------------
$newRows = @(
    @{ Row = 22; D = 44428; M = 100; N = 20000; O = 21000; P = 20500; S = 1025 },
    @{ Row = 23; D = 44442; M = 140; N = 20000; O = 21000; P = 20500; S = 1025 },
    @{ Row = 24; D = 44435; M = 260; N = 20000; O = 22000; P = 21115; S = 1056 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = 8
    $ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"

    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Range("D" + $row).NumberFormat = $ws.Range("D21").NumberFormat

    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value = 100108007
    $ws.Cells.Item($row, 10).Value = "Coco"
    $ws.Cells.Item($row, 11).Value = "Sin especificar"
    $ws.Cells.Item($row, 12).Value = "Primera"

    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P

    $ws.Cells.Item($row, 17).Value = "$/malla 20 unidades"
    $ws.Cells.Item($row, 18).Value = "Perú"

    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = 20
}
